$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    $rng = $ws.Range($cellAddr)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

$ws.Range('D2').Value = '64.303.13'
$ws.Range('E2').Value = '  +0.45%  '
$ws.Range('D3').Value = '3.329.50'
$ws.Range('E3').Value = '  +0.05%  '
Set-TextValue 'D4' '1.00'
$ws.Range('E4').Value = '  +0.00%  '
Set-TextValue 'D5' '553.55'
$ws.Range('E5').Value = '  +0.38%  '
Set-TextValue 'D6' '173.39'
$ws.Range('E6').Value = '  +0.54%  '
Set-TextValue 'D7' '0.620'
$ws.Range('E7').Value = '  +1.15%  '
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('D9').Value = '3.321.32'
$ws.Range('E9').Value = '  -0.01%  '
$ws.Range('E10').Value = '  +5.29%  '
Set-TextValue 'D11' '0.632'
$ws.Range('E11').Value = '  +1.79%  '
Set-TextValue 'D12' '53.37'
$ws.Range('E12').Value = '  +0.57%  '
Set-TextValue 'D13' '0.0000279'
$ws.Range('E13').Value = '  +2.67%  '
Set-TextValue 'D14' '9.10'
$ws.Range('E14').Value = '  +1.19%  '
$ws.Range('D15').Value = '3.852.96'
$ws.Range('E15').Value = '  -0.22%  '
$ws.Range('E16').Value = '  +3.04%  '
Set-TextValue 'D17' '18.11'
$ws.Range('E17').Value = '  -0.68%  '
$ws.Range('D18').Value = '3.349.66'
$ws.Range('E18').Value = '  +0.95%  '
$ws.Range('D19').Value = '64.138.78'
$ws.Range('E19').Value = '  +0.34%  '
Set-TextValue 'D20' '11.71'
$ws.Range('E20').Value = '  -0.14%  '
Set-TextValue 'D21' '0.986'
$ws.Range('E21').Value = '  +1.58%  '
Set-TextValue 'D22' '453.83'
$ws.Range('E22').Value = '  +6.28%  '
Set-TextValue 'D23' '5.10'
$ws.Range('E23').Value = '  +9.49%  '
$ws.Range('E24').Value = '  -0.29%  '
$ws.Range('B25').Value = 'InternetComputer(DFINITY)'
$ws.Range('C25').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue 'D25' '14.07'
$ws.Range('E25').Value = '  +6.07%  '
$ws.Range('B26').Value = 'Litecoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue 'D26' '87.30'
$ws.Range('E26').Value = '  +3.89%  '
Set-TextValue 'D27' '2.87'
$ws.Range('E27').Value = '  +2.29%  '
$ws.Range('E28').Value = '  -0.10%  '
Set-TextValue 'D29' '30.95'
$ws.Range('E29').Value = '  +4.45%  '
Set-TextValue 'D30' '8.58'
$ws.Range('E30').Value = '  +0.01%  '
$ws.Range('E31').Value = '  -2.35%  '
Set-TextValue 'D32' '11.41'
$ws.Range('E32').Value = '  +0.30%  '
Set-TextValue 'D33' '61.64'
$ws.Range('E33').Value = '  +6.11%  '
Set-TextValue 'D34' '565.85'
$ws.Range('E34').Value = '  -4.80%  '
$ws.Range('E35').Value = '  +0.12%  '
$ws.Range('E36').Value = '  -0.06%  '
Set-TextValue 'D37' '0.141'
$ws.Range('E37').Value = '  -1.35%  '
Set-TextValue 'D38' '3.50'
$ws.Range('E38').Value = '  +1.55%  '
Set-TextValue 'D39' '35.30'
$ws.Range('E39').Value = '  +0.12%  '
$ws.Range('E40').Value = '  +0.34%  '
$ws.Range('D41').Value = '0.0₃0728'
$ws.Range('E41').Value = '  -2.68%  '
$ws.Range('D42').Value = '3.058.77'
$ws.Range('E42').Value = '  -0.96%  '
$ws.Range('E43').Value = '  +2.79%  '
Set-TextValue 'D44' '2.75'
$ws.Range('E44').Value = '  -0.98%  '
Set-TextValue 'D45' '3.19'
$ws.Range('E45').Value = '  +0.25%  '
Set-TextValue 'D46' '2.45'
$ws.Range('E46').Value = '  +0.70%  '
$ws.Range('E47').Value = '  +3.48%  '
$ws.Range('E48').Value = '  -0.02%  '
Set-TextValue 'D49' '141.72'
$ws.Range('E49').Value = '  +6.84%  '
$ws.Range('E50').Value = '  -3.38%  '
Set-TextValue 'D51' '8.16'
$ws.Range('E51').Value = '  +0.72%  '
